# Auto-generated Excel COM-interop edit script
# Applies market-data refresh updates to the Yojimbo_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$ws.Range("H40").Value = 51500
$ws.Range("I40").Value = 100000
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 100000
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -99825
$ws.Range("N40").Value = -3350

# Row 100 (ALC)
$ws.Range("H100").Value = 2126.8462
$ws.Range("I100").Value = 1604.4546
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 1604.4546
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -1063.4546
$ws.Range("N100").Value = -6082

# Row 137 (ALC)
$ws.Range("H137").Value = 2019.0476
$ws.Range("I137").Value = 2024.8269
$ws.Range("J137").Value = 1991.7273
$ws.Range("K137").Value = 6074.4807
$ws.Range("L137").Value = 5975.1819
$ws.Range("M137").Value = -3524.4807
$ws.Range("N137").Value = -11075.1819

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 6023.25
$ws.Range("I61").Value = 6023.25
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6023.25
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5811.25
$ws.Range("N61").ClearContents()

# Row 76 (ARM)
$ws.Range("H76").Value = 20910
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 20910
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 20910
$ws.Range("N76").Value = -21586

# Row 79 (ARM)
$ws.Range("H79").Value = 20910
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 20910
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 20910
$ws.Range("N79").Value = -23250

# Row 97 (ARM)
$ws.Range("H97").Value = 1898.826
$ws.Range("I97").Value = 934.2632
$ws.Range("J97").Value = 6480.5
$ws.Range("K97").Value = 934.2632
$ws.Range("L97").Value = 6480.5
$ws.Range("M97").Value = -438.2632
$ws.Range("N97").Value = -7472.5

# Row 136 (ARM)
$ws.Range("H136").Value = 6023.25
$ws.Range("I136").Value = 6023.25
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 18069.75
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15519.75
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 40 (BSM)
$ws.Range("H40").Value = 38298.668
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 38298.668
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 38298.668
$ws.Range("N40").Value = -38828.668

# Row 76 (BSM)
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79 (BSM)
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 86 (BSM)
$ws.Range("H86").Value = 2062.6924
$ws.Range("I86").Value = 2279.7
$ws.Range("J86").Value = 1339.3334
$ws.Range("K86").Value = 2279.7
$ws.Range("L86").Value = 1339.3334
$ws.Range("M86").Value = -1156.7
$ws.Range("N86").Value = -3585.3334

# Row 89 (BSM)
$ws.Range("H89").Value = 2062.6924
$ws.Range("I89").Value = 2279.7
$ws.Range("J89").Value = 1339.3334
$ws.Range("K89").Value = 11398.5
$ws.Range("L89").Value = 6696.666999999999
$ws.Range("M89").Value = -5782.5
$ws.Range("N89").Value = -17928.667

# Row 94 (BSM)
$ws.Range("H94").Value = 1138.5172
$ws.Range("I94").Value = 822.0357
$ws.Range("J94").Value = 10000
$ws.Range("K94").Value = 822.0357
$ws.Range("L94").Value = 10000
$ws.Range("M94").Value = -371.0357
$ws.Range("N94").Value = -10902

# Row 99 (BSM)
$ws.Range("H99").Value = 1180.7391
$ws.Range("I99").Value = 1122.8823
$ws.Range("J99").Value = 1344.6666
$ws.Range("K99").Value = 1122.8823
$ws.Range("L99").Value = 1344.6666
$ws.Range("M99").Value = 375.1177
$ws.Range("N99").Value = -4340.6666

# Row 134 (BSM)
$ws.Range("H134").Value = 4360.5713
$ws.Range("I134").Value = 5793
$ws.Range("J134").Value = 1936.4615
$ws.Range("K134").Value = 17379
$ws.Range("L134").Value = 5809.3845
$ws.Range("M134").Value = -14844
$ws.Range("N134").Value = -10879.3845

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Range("H4").Value = 48000.8
$ws.Range("I4").Value = 33333.332
$ws.Range("J4").Value = 70002
$ws.Range("K4").Value = 33333.332
$ws.Range("L4").Value = 70002
$ws.Range("M4").Value = -33221.332
$ws.Range("N4").Value = -70226

# Row 31 (CRP)
$ws.Range("H31").Value = 24081.408
$ws.Range("I31").Value = 37654.4
$ws.Range("J31").Value = 2650.3684
$ws.Range("K31").Value = 37654.4
$ws.Range("L31").Value = 2650.3684
$ws.Range("M31").Value = -37359.4
$ws.Range("N31").Value = -3240.3684

# Row 34 (CRP)
$ws.Range("H34").Value = 24081.408
$ws.Range("I34").Value = 37654.4
$ws.Range("J34").Value = 2650.3684
$ws.Range("K34").Value = 37654.4
$ws.Range("L34").Value = 2650.3684
$ws.Range("M34").Value = -37452.4
$ws.Range("N34").Value = -3054.3684

# Row 58 (CRP)
$ws.Range("H58").Value = 1265.6052
$ws.Range("I58").Value = 1316.5
$ws.Range("J58").Value = 1074.75
$ws.Range("K58").Value = 1316.5
$ws.Range("L58").Value = 1074.75
$ws.Range("M58").Value = -1113.5
$ws.Range("N58").Value = -1480.75

# Row 86 (CRP)
$ws.Range("H86").Value = 2543.1
$ws.Range("I86").Value = 2045.6
$ws.Range("J86").Value = 3040.6
$ws.Range("K86").Value = 2045.6
$ws.Range("L86").Value = 3040.6
$ws.Range("M86").Value = -922.5999999999999
$ws.Range("N86").Value = -5286.6

# Row 89 (CRP)
$ws.Range("H89").Value = 2543.1
$ws.Range("I89").Value = 2045.6
$ws.Range("J89").Value = 3040.6
$ws.Range("K89").Value = 10228
$ws.Range("L89").Value = 15203
$ws.Range("M89").Value = -4612
$ws.Range("N89").Value = -26435

# Row 95 (CRP)
$ws.Range("H95").Value = 17653.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17653.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17653.5
$ws.Range("N95").Value = -23145.5

# Row 96 (CRP)
$ws.Range("H96").Value = 13000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 13000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 13000
$ws.Range("N96").Value = -18492

# Row 136 (CRP)
$ws.Range("H136").Value = 1265.6052
$ws.Range("I136").Value = 1316.5
$ws.Range("J136").Value = 1074.75
$ws.Range("K136").Value = 3949.5
$ws.Range("L136").Value = 3224.25
$ws.Range("M136").Value = -1399.5
$ws.Range("N136").Value = -8324.25

$ws = $wb.Worksheets.Item("GSM")
# Row 95 (GSM)
$ws.Range("H95").Value = 16292.143
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 16292.143
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 16292.143
$ws.Range("N95").Value = -21784.143

$ws = $wb.Worksheets.Item("LTW")
# Row 17 (LTW)
$ws.Range("H17").Value = 1350166.6
$ws.Range("I17").Value = 2000250
$ws.Range("J17").Value = 50000
$ws.Range("K17").Value = 2000250
$ws.Range("L17").Value = 50000
$ws.Range("M17").Value = -2000080
$ws.Range("N17").Value = -50340

# Row 136 (LTW)
$ws.Range("H136").Value = 3488.2727
$ws.Range("I136").Value = 3018.9375
$ws.Range("J136").Value = 3930
$ws.Range("K136").Value = 9056.8125
$ws.Range("L136").Value = 11790
$ws.Range("M136").Value = -6506.8125
$ws.Range("N136").Value = -16890

$ws = $wb.Worksheets.Item("WVR")
# Row 55 (WVR)
$ws.Range("H55").Value = 14577.714
$ws.Range("I55").Value = 2136
$ws.Range("J55").Value = 31166.666
$ws.Range("K55").Value = 2136
$ws.Range("L55").Value = 31166.666
$ws.Range("M55").Value = -1859
$ws.Range("N55").Value = -31720.666

# Row 136 (WVR)
$ws.Range("H136").Value = 4334.778
$ws.Range("I136").Value = 4383.7354
$ws.Range("J136").Value = 3502.5
$ws.Range("K136").Value = 13151.2062
$ws.Range("L136").Value = 10507.5
$ws.Range("M136").Value = -10601.2062
$ws.Range("N136").Value = -15607.5
